# Updates cryptos list (prices / 1h volume %) and swaps the Chainlink/WrappedEther
# and Filecoin/Monero row pairs, per the Jan 18 2024 GitHub Actions refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain (non-formula) value into a cell while keeping it a
# text/string cell. Digit-and-dot-only strings (e.g. "1.00", "311.28") would
# otherwise be auto-coerced by Excel into numeric cells, which both changes
# the stored type and can lose formatting like trailing zeros. Forcing the
# cell to Text first keeps the literal string, and resetting the Style back
# to Normal afterwards avoids leaving a stray custom number format behind.
function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-TextCell "D2" "42.825.44"
Set-Cell "E2" "  +0.36%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.540.43"
Set-Cell "E3" "  -0.15%  "

# Row 4 - TetherUSD
Set-TextCell "D4" "1.00"
Set-Cell "E4" "  -0.06%  "

# Row 5 - BNB
Set-TextCell "D5" "311.28"
Set-Cell "E5" "  +0.76%  "

# Row 6 - Solana
Set-TextCell "D6" "100.81"
Set-Cell "E6" "  +3.50%  "

# Row 7 - XRP
Set-Cell "E7" "  -0.71%  "

# Row 8 - USDC
Set-Cell "E8" "  +0.04%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.524"
Set-Cell "E9" "  -0.85%  "

# Row 10 - Avalanche
Set-TextCell "D10" "35.83"
Set-Cell "E10" "  +1.12%  "

# Row 11 - Dogecoin
Set-TextCell "D11" "0.0806"

# Row 12 - Polkadot
Set-Cell "E12" "  -0.70%  "

# Row 13 - TRON
Set-Cell "E13" "  +1.72%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell "D14" "2.928.76"
Set-Cell "E14" "  -0.21%  "

# Rows 15 & 16 swap places: Chainlink <-> WrappedEther
Set-Cell "B15" "WrappedEther"
Set-Cell "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D15" "2.566.41"
Set-Cell "E15" "  +0.28%  "

Set-Cell "B16" "Chainlink"
Set-Cell "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D16" "15.37"
Set-Cell "E16" "  -2.14%  "

# Row 17 - Polygon
Set-Cell "E17" "  -1.75%  "

# Row 18 - WrappedBTC
Set-TextCell "D18" "42.824.22"
Set-Cell "E18" "  +0.27%  "

# Row 19 - Uniswap
Set-TextCell "D19" "6.74"
Set-Cell "E19" "  +0.36%  "

# Row 20 - InternetComputer(DFINITY)
Set-TextCell "D20" "12.39"
Set-Cell "E20" "  +0.23%  "

# Row 21 - ShibaInu
Set-TextCell "D21" "0.0₃0954"
Set-Cell "E21" "  -0.09%  "

# Row 22 - Litecoin
Set-Cell "E22" "  +1.44%  "

# Row 23 - BitcoinCash
Set-TextCell "D23" "244.00"
Set-Cell "E23" "  -1.18%  "

# Row 24 - PancakeSwap
Set-TextCell "D24" "2.88"
Set-Cell "E24" "  -0.99%  "

# Row 25 - ImmutableX
Set-Cell "E25" "  -0.18%  "

# Row 26 - Dai
Set-Cell "E26" "  +0.02%  "

# Row 27 - EthereumClassic
Set-TextCell "D27" "25.53"
Set-Cell "E27" "  -3.69%  "

# Row 28 - Toncoin
Set-Cell "E28" "  -1.22%  "

# Row 29 - Cosmos
Set-TextCell "D29" "10.18"
Set-Cell "E29" "  +0.52%  "

# Row 30 - InjectiveProtocol
Set-TextCell "D30" "38.66"
Set-Cell "E30" "  -3.82%  "

# Rows 31 & 32 swap places: Filecoin <-> Monero
Set-Cell "B31" "Monero"
Set-Cell "C31" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D31" "158.72"
Set-Cell "E31" "  +0.82%  "

Set-Cell "B32" "Filecoin"
Set-Cell "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D32" "5.87"
Set-Cell "E32" "  +2.70%  "

# Row 33 - ApeXProtocol
Set-TextCell "D33" "2.77"
Set-Cell "E33" "  +7.41%  "

# Row 34 - WEMIXToken
Set-Cell "E34" "  +2.32%  "

# Row 35 - Hedera
Set-Cell "E35" "  +0.20%  "

# Row 36 - Celestia
Set-TextCell "D36" "18.24"
Set-Cell "E36" "  -0.79%  "

# Row 37 - LidoDAOToken
Set-TextCell "D37" "3.16"
Set-Cell "E37" "  -3.39%  "

# Row 38 - ARBITRUM
Set-TextCell "D38" "1.97"
Set-Cell "E38" "  -4.42%  "

# Row 39 - Kaspa
Set-Cell "E39" "  +0.42%  "

# Row 40 - Stellar
Set-Cell "E40" "  +0.24%  "

# Row 41 - RenderToken
Set-TextCell "D41" "4.16"
Set-Cell "E41" "  +3.21%  "

# Row 42 - EnergySwap
Set-TextCell "D42" "21.91"
Set-Cell "E42" "  -1.66%  "

# Row 44 - NEARProtocol
Set-TextCell "D44" "3.32"
Set-Cell "E44" "  +4.08%  "

# Row 45 - VeChain
Set-Cell "E45" "  +0.47%  "

# Row 46 - Maker
Set-TextCell "D46" "2.003.34"
Set-Cell "E46" "  +0.68%  "

# Row 47 - FraxShare
Set-TextCell "D47" "9.06"
Set-Cell "E47" "  +0.11%  "

# Row 48 - RocketPoolETH
Set-TextCell "D48" "2.780.28"

# Row 49 - Algorand
Set-Cell "E49" "  +0.55%  "

# Row 50 - BitcoinSV
Set-TextCell "D50" "80.32"
Set-Cell "E50" "  -0.35%  "

# Row 51 - ordi
Set-TextCell "D51" "72.64"
Set-Cell "E51" "  -0.62%  "
